# Update the cryptos price table (Sheet1) with refreshed price / volume
# figures, matching the latest GitHub Actions data refresh.
# Column D values are forced to Text format before assignment so that
# numeric-looking strings (e.g. "0.9975", "106.00", "26.615.64") are
# preserved exactly as text instead of being coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.615.64"
$ws.Cells.Item(2, 5).Value = "  -0.07%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.720.28"
$ws.Cells.Item(3, 5).Value = "  -1.14%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9975"
$ws.Cells.Item(4, 5).Value = "  -0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "241.16"
$ws.Cells.Item(5, 5).Value = "  -2.04%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9980"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4938"
$ws.Cells.Item(7, 5).Value = "  -0.32%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2608"
$ws.Cells.Item(8, 5).Value = "  -2.59%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06213"
$ws.Cells.Item(9, 5).Value = "  -0.87%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.726.38"
$ws.Cells.Item(10, 5).Value = "  -0.77%  "
$ws.Cells.Item(11, 2).Value = "Solana"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "15.81"
$ws.Cells.Item(11, 5).Value = "  +0.26%  "
$ws.Cells.Item(12, 2).Value = "TRON"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.07008"
$ws.Cells.Item(12, 5).Value = "  -0.53%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.6082"
$ws.Cells.Item(13, 5).Value = "  -1.00%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.497"
$ws.Cells.Item(14, 5).Value = "  -2.34%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "76.92"
$ws.Cells.Item(15, 5).Value = "  -1.46%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.9982"
$ws.Cells.Item(16, 5).Value = "  -0.16%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "26.455.93"
$ws.Cells.Item(17, 5).Value = "  -0.72%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.9974"
$ws.Cells.Item(18, 5).Value = "  -0.24%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.000007165"
$ws.Cells.Item(19, 5).Value = "  -1.20%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.37"
$ws.Cells.Item(20, 5).Value = "  -1.69%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.943.36"
$ws.Cells.Item(21, 5).Value = "  -1.18%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.419"
$ws.Cells.Item(22, 5).Value = "  -3.00%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.526"
$ws.Cells.Item(23, 5).Value = "  -2.39%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.092"
$ws.Cells.Item(24, 5).Value = "  -3.87%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "137.88"
$ws.Cells.Item(25, 5).Value = "  -0.82%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "15.31"
$ws.Cells.Item(26, 5).Value = "  -0.51%  "
$ws.Cells.Item(27, 5).Value = "  -0.77%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.737"
$ws.Cells.Item(28, 5).Value = "  -1.24%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "106.00"
$ws.Cells.Item(29, 5).Value = "  -1.20%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.927"
$ws.Cells.Item(30, 5).Value = "  -2.30%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.07962"
$ws.Cells.Item(31, 5).Value = "  -0.80%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.655"
$ws.Cells.Item(32, 5).Value = "  -1.79%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04501"
$ws.Cells.Item(33, 5).Value = "  -2.36%  "
$ws.Cells.Item(34, 2).Value = "Frax"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.9972"
$ws.Cells.Item(34, 5).Value = "  -0.19%  "
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.611"
$ws.Cells.Item(35, 5).Value = "  -0.10%  "
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.0000"
$ws.Cells.Item(36, 5).Value = "  -1.73%  "
$ws.Cells.Item(37, 2).Value = "ImmutableX"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.6243"
$ws.Cells.Item(37, 5).Value = "  -2.09%  "
$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.9369"
$ws.Cells.Item(38, 5).Value = "  +2.89%  "
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.997"
$ws.Cells.Item(39, 5).Value = "  -2.62%  "
$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.417"
$ws.Cells.Item(40, 5).Value = "  -0.51%  "
$ws.Cells.Item(41, 2).Value = "PaxDollar"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.9978"
$ws.Cells.Item(41, 5).Value = "  -0.49%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.01511"
$ws.Cells.Item(42, 5).Value = "  +0.34%  "
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.516"
$ws.Cells.Item(43, 5).Value = "  +0.98%  "
$ws.Cells.Item(44, 2).Value = "Quant"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "99.26"
$ws.Cells.Item(44, 5).Value = "  -2.55%  "
$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.3837"
$ws.Cells.Item(45, 5).Value = "  -2.56%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "6.917"
$ws.Cells.Item(46, 5).Value = "  +0.89%  "
$ws.Cells.Item(47, 2).Value = "Algorand"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.1159"
$ws.Cells.Item(47, 5).Value = "  -1.61%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.05374"
$ws.Cells.Item(48, 5).Value = "  -0.17%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.742"
$ws.Cells.Item(49, 5).Value = "  -0.78%  "
$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "30.12"
$ws.Cells.Item(50, 5).Value = "  -1.94%  "
$ws.Cells.Item(51, 2).Value = "Aave"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "51.48"
$ws.Cells.Item(51, 5).Value = "  -0.59%  "
